# Update TPM-derived NATMI metrics for Spon2-Itga4 ligand-receptor pairs
# (recalculated with new TPM values; commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1079986666666667
$ws.Range("H2").Value = 0.323996
$ws.Range("I2").Value = 0.004187739561209694
$ws.Range("J2").Value = 0.004187739561209694
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.002514856952
$ws.Range("R2").Value = 0.022633712568
$ws.Range("S2").Value = 0.00003899083060245817
$ws.Range("T2").Value = 0.00003899083060245818

# Row 3
$ws.Range("G3").Value = 0.1079986666666667
$ws.Range("H3").Value = 0.323996
$ws.Range("I3").Value = 0.004187739561209694
$ws.Range("J3").Value = 0.004187739561209694
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.01435511077422222
$ws.Range("R3").Value = 0.129195996968
$ws.Range("S3").Value = 0.0002225644253968767
$ws.Range("T3").Value = 0.0002225644253968767

# Row 4
$ws.Range("G4").Value = 0.1079986666666667
$ws.Range("H4").Value = 0.323996
$ws.Range("I4").Value = 0.004187739561209694
$ws.Range("J4").Value = 0.004187739561209694
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 0.2532336896195556
$ws.Range("R4").Value = 2.279103206576
$ws.Range("S4").Value = 0.00392618430521036
$ws.Range("T4").Value = 0.00392618430521036

# Row 5
$ws.Range("I5").Value = 0.9687110856121154
$ws.Range("J5").Value = 0.9687110856121155
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.581738614
$ws.Range("R5").Value = 5.235647526
$ws.Range("S5").Value = 0.009019388452827912
$ws.Range("T5").Value = 0.009019388452827916

# Row 6
$ws.Range("I6").Value = 0.9687110856121154
$ws.Range("J6").Value = 0.9687110856121155
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.05148377137439881
$ws.Range("T6").Value = 0.05148377137439882

# Row 7
$ws.Range("I7").Value = 0.9687110856121154
$ws.Range("J7").Value = 0.9687110856121155
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 58.57820879244446
$ws.Range("R7").Value = 527.2038791320001
$ws.Range("S7").Value = 0.9082079257848887
$ws.Range("T7").Value = 0.9082079257848888

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6989190000000001
$ws.Range("H8").Value = 2.096757
$ws.Range("I8").Value = 0.02710117482667488
$ws.Range("J8").Value = 0.02710117482667489
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.016275027834
$ws.Range("R8").Value = 0.146475250506
$ws.Range("S8").Value = 0.0002523311923650859
$ws.Range("T8").Value = 0.0002523311923650861

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6989190000000001
$ws.Range("H9").Value = 2.096757
$ws.Range("I9").Value = 0.02710117482667488
$ws.Range("J9").Value = 0.02710117482667489
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 0.09289984753399999
$ws.Range("R9").Value = 0.836098627806
$ws.Range("S9").Value = 0.001440337278552448
$ws.Range("T9").Value = 0.001440337278552449

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6989190000000001
$ws.Range("H10").Value = 2.096757
$ws.Range("I10").Value = 0.02710117482667488
$ws.Range("J10").Value = 0.02710117482667489
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 1.638815020388001
$ws.Range("R10").Value = 14.749335183492
$ws.Range("S10").Value = 0.02540850635575735
$ws.Range("T10").Value = 0.02540850635575735

